# Add masking password feature
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DATA")
$ws.Activate()

# Mask the plaintext password values with a base64-encoded representation
$ws.Range("E2").Value = "YWRtaW4xMjM="
$ws.Range("E3").Value = "YWRtaW4xMjM="

# Flip the "execute" flag for the amazonTest row from yes to no
$ws.Range("B7").Value = "no"

# Update the active selection on the DATA sheet
$ws.Range("F9").Select()
